$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sprint 2")
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Sprint 3"

$ws2.Range("B2:W48").Copy()
$ws3.Range("B2:W48").PasteSpecial(-4122)
Write-Output "pasted"
